$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new tracking columns: frame_w / frame_h (the source frame
# resolution parsed out of each row's Stream filename, e.g.
# ".\sources\RaceHorses_416x240_30.yuv" -> 416 x 240).
$ws.Range("I1").Value = "frame_w"
$ws.Range("J1").Value = "frame_h"

# Match the column widths used for the rest of the sheet (values are
# fed in "ColumnWidth" units; this engine stores width = ColumnWidth +
# 5/6 in the saved XML, so back the offset out to land on 9 / 8.57).
$ws.Columns.Item(9).ColumnWidth = 9 - 5 / 6
$ws.Columns.Item(10).ColumnWidth = 8.5703125 - 5 / 6

# frame_w / frame_h per example row, derived from the Stream column.
$frameDims = @{
    2  = @(416, 240)   # RaceHorses_416x240_30.yuv
    3  = @(416, 240)   # RaceHorses_416x240_30.yuv
    4  = @(432, 240)   # VQ_sample_432x240.yuv
    5  = @(432, 240)   # VQ_sample_432x240.yuv
    6  = @(416, 240)   # RaceHorses_416x240_30.yuv
    7  = @(432, 240)   # VQ_sample_432x240.yuv
    8  = @(432, 240)   # VQ_sample_432x240.yuv
    9  = @(416, 240)   # BasketballPass_416x240_50.yuv
    10 = @(416, 240)   # BasketballPass_416x240_50.yuv
    11 = @(416, 240)   # BasketballPass_416x240_50.yuv
    12 = @(416, 240)   # BasketballPass_416x240_50.yuv
    13 = @(432, 240)   # VQ_sample_432x240.yuv
    14 = @(432, 240)   # VQ_sample_432x240.yuv
    15 = @(416, 240)   # RaceHorses_416x240_30.yuv
    16 = @(416, 240)   # RaceHorses_416x240_30.yuv
    17 = @(416, 240)   # RaceHorses_416x240_30.yuv
    18 = @(416, 240)   # RaceHorses_416x240_30.yuv
    19 = @(416, 240)   # BasketballPass_416x240_50.yuv
    20 = @(416, 240)   # BasketballPass_416x240_50.yuv
    21 = @(432, 240)   # VQ_sample_432x240.yuv
    22 = @(432, 240)   # VQ_sample_432x240.yuv
    23 = @(432, 240)   # VQ_sample_432x240.yuv
    24 = @(416, 240)   # BasketballPass_416x240_50.yuv
    25 = @(416, 240)   # BasketballPass_416x240_50.yuv
    26 = @(416, 240)   # RaceHorses_416x240_30.yuv
    27 = @(416, 240)   # RaceHorses_416x240_30.yuv
    28 = @(416, 240)   # BasketballPass_416x240_50.yuv
}

foreach ($r in 2..28) {
    $dims = $frameDims[$r]
    $ws.Range("I$r").Value = $dims[0]
    $ws.Range("J$r").Value = $dims[1]
}
